# Full Planning Permission - Full: add "Trees on adjacent land" row under
# "Trees and hedges information" and fill in the first sub-item
# ("Trees on site") which was previously blank. This pushes the existing
# "Vehicle parking", "Waste storage and collection" blocks (and everything
# below) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 207; this shifts rows 207-215 down to 208-216,
# copying formatting from the row above (keeps style s="2") and keeps
# merged-cell ranges for A208:A212 / B208:B212 / A213:A216 / B213:B216
# (shifted from the old A207:A211 / B207:B211 / A212:A215 / B212:B215),
# while the A206/B206 single-cell merges grow to A206:A207 / B206:B207
# automatically because the insertion happens inside that merged block.
$ws.Rows.Item(207).Insert()

# Row 206: "Trees and hedges information" / first sub-item -> "Trees on site"
$ws.Range("C206").Value = "Trees on site"
$ws.Range("G206").Value = "Whether trees or hedges are present on the proposed development site"
$ws.Range("H206").Value = "boolean"
$ws.Range("I206").Value = "MUST"

# Row 207 (new): second sub-item -> "Trees on adjacent land"
$ws.Range("C207").Value = "Trees on adjacent land"
$ws.Range("G207").Value = "Whether trees or hedges on land adjacent to the proposed development site could influence the development or might be important as part of the local landscape character"
$ws.Range("H207").Value = "boolean"
$ws.Range("I207").Value = "MUST"
